$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated correlation matrix values (rows 2-12, columns B-L)
$values = @(
    @(1, -0.1625239993632729, 0.03034932294482378, 0.1252628789284959, -0.1605040288970766, -0.09995119463291788, -0.06084151847845858, 0.2427483138624487, -0.1596361618078241, -0.05455505983836755, 0.08587227291144571),
    @(-0.1625239993632729, 1, -0.1652764414883613, 0.07155001752140763, 0.08318169954040049, -0.05450870907588182, -0.1480642479747007, -0.1067363692665227, 0.5195342058528954, 0.1229351664332261, -0.1535429869349732),
    @(0.03034932294482378, -0.1652764414883613, 1, -0.107106122654923, -0.08532197024002931, 0.1379330680931256, 0.2123100299950188, 0.04313813038796596, -0.1524330778905541, -0.2751239234336399, 0.1797298776755977),
    @(0.1252628789284959, 0.07155001752140763, -0.107106122654923, 1, -0.0593406423198977, -0.08777297298871735, -0.2414084925065061, 0.04082888771833498, 0.05050445067887077, 0.007827052947800506, -0.02273158211910532),
    @(-0.1605040288970766, 0.08318169954040049, -0.08532197024002931, -0.0593406423198977, 1, -0.02128189587790422, 0.002313369359746511, -0.1939715910740801, 0.09396031763861029, 0.04488902142332702, -0.03046516221179131),
    @(-0.09995119463291788, -0.05450870907588182, 0.1379330680931256, -0.08777297298871735, -0.02128189587790422, 1, 0.1581061620943013, 0.02379014056509399, -0.04030548155315947, -0.1577074545228274, 0.1393915103169893),
    @(-0.06084151847845858, -0.1480642479747007, 0.2123100299950188, -0.2414084925065061, 0.002313369359746511, 0.1581061620943013, 1, 0.04775404509345622, -0.1035072310202461, -0.1820342848792459, 0.2786310469845365),
    @(0.2427483138624487, -0.1067363692665227, 0.04313813038796596, 0.04082888771833498, -0.1939715910740801, 0.02379014056509399, 0.04775404509345622, 1, -0.1088731746243381, -0.07072663486603732, 0.1271916454228323),
    @(-0.1596361618078241, 0.5195342058528954, -0.1524330778905541, 0.05050445067887077, 0.09396031763861029, -0.04030548155315947, -0.1035072310202461, -0.1088731746243381, 1, 0.1051393112083464, -0.1245972999285731),
    @(-0.05455505983836755, 0.1229351664332261, -0.2751239234336399, 0.007827052947800506, 0.04488902142332702, -0.1577074545228274, -0.1820342848792459, -0.07072663486603732, 0.1051393112083464, 1, -0.3062534593909333),
    @(0.08587227291144571, -0.1535429869349732, 0.1797298776755977, -0.02273158211910532, -0.03046516221179131, 0.1393915103169893, 0.2786310469845365, 0.1271916454228323, -0.1245972999285731, -0.3062534593909333, 1)
)

for ($i = 0; $i -lt 11; $i++) {
    for ($j = 0; $j -lt 11; $j++) {
        $row = $i + 2
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$i][$j]
    }
}
